$wb = $excel.ActiveWorkbook

# --- Sheet "All Orders": update row 17 (Order ID 3) ---
$wsOrders = $wb.Worksheets.Item("All Orders")
$wsOrders.Range("H17").Value = "CANCELLED"
$wsOrders.Range("M17").Value = "test order"

# --- Sheet "Daily Summary": update row 4 (2026-01-13) totals ---
$wsSummary = $wb.Worksheets.Item("Daily Summary")
$wsSummary.Range("D4").Value = 4
$wsSummary.Range("E4").Value = 170
$wsSummary.Range("G4").Value = 170
